$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ahuMdl")
$ws2 = $wb.Worksheets.Item("faults")

# --- sheet "faults": update header text ---
$ws2.Range("C1").Value = "AHU Health Index (%)"

# --- sheet "faults": row 3 (AHU2) updates ---
$ws2.Range("C3").Value = 33.33333333333334
$ws2.Range("G3").Value = "Low outdoor air"

# --- sheet "ahuMdl": updated numeric model output values ---
$ws1.Range("B2").Value = -15.90307908750653
$ws1.Range("C2").Value = -2.511870820920659
$ws1.Range("D2").Value = 0.4928514705438788

$ws1.Range("B3").Value = -12.25942311130128
$ws1.Range("C3").Value = 8.550331952986371
$ws1.Range("D3").Value = 0.252985159787625
